# Update the "Buying Opportunity" (column B) and "support Zone" (column C)
# ticker lists on Sheet1, then drop the now-unused trailing rows (37-40) so
# the sheet's dimension shrinks from A1:F40 to A1:F36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B / Column C updates for rows 2-12 (both columns change) ---
$ws.Range("B2").Value  = "NSE:AARTIPHARM"
$ws.Range("C2").Value  = "NSE:ASPINWALL"

$ws.Range("B3").Value  = "NSE:AKSHARCHEM"
$ws.Range("C3").Value  = "NSE:AVTNPL"

$ws.Range("B4").Value  = "NSE:ARVEE"
$ws.Range("C4").Value  = "NSE:CAMPUS"

$ws.Range("B5").Value  = "NSE:AXISHCETF"
$ws.Range("C5").Value  = "NSE:ENIL"

$ws.Range("B6").Value  = "NSE:BANKINDIA"
$ws.Range("C6").Value  = "NSE:INOXWIND"

$ws.Range("B7").Value  = "NSE:BANSWRAS"
$ws.Range("C7").Value  = "NSE:KRISHANA"

$ws.Range("B8").Value  = "NSE:BDL"
$ws.Range("C8").Value  = "NSE:MADHAV"

$ws.Range("B9").Value  = "NSE:BLUEDART"
$ws.Range("C9").Value  = "NSE:MCLEODRUSS"

$ws.Range("B10").Value = "NSE:CASTROLIND"
$ws.Range("C10").Value = "NSE:PODDARMENT"

$ws.Range("B11").Value = "NSE:CMSINFO"
$ws.Range("C11").Value = "NSE:RHL"

$ws.Range("B12").Value = "NSE:DEVIT"
$ws.Range("C12").Value = "NSE:RITES"

# --- Column B only updates for rows 13-36 (column C stays blank) ---
$ws.Range("B13").Value = "NSE:EBBETF0431"
$ws.Range("B14").Value = "NSE:FDC"
$ws.Range("B15").Value = "NSE:FINCABLES"
$ws.Range("B16").Value = "NSE:GOLDTECH"
$ws.Range("B17").Value = "NSE:GREENPOWER"
$ws.Range("B18").Value = "NSE:HDFCLIQUID"
$ws.Range("B19").Value = "NSE:HEALTHY"
$ws.Range("B20").Value = "NSE:HGS"
$ws.Range("B21").Value = "NSE:HPIL"
$ws.Range("B22").Value = "NSE:INDNIPPON"
$ws.Range("B23").Value = "NSE:INTLCONV"
$ws.Range("B24").Value = "NSE:IRMENERGY"
$ws.Range("B25").Value = "NSE:IVP"
$ws.Range("B26").Value = "NSE:JKLAKSHMI"
$ws.Range("B27").Value = "NSE:LORDSCHLO"
$ws.Range("B28").Value = "NSE:MAFANG"
$ws.Range("B29").Value = "NSE:MAPMYINDIA"
$ws.Range("B30").Value = "NSE:MONQ50"
$ws.Range("B31").Value = "NSE:PATINTLOG"
$ws.Range("B32").Value = "NSE:PSPPROJECT"
$ws.Range("B33").Value = "NSE:PTL"
$ws.Range("B34").Value = "NSE:RELAXO"
$ws.Range("B35").Value = "NSE:RELIGARE"
$ws.Range("B36").Value = "NSE:RUSTOMJEE"

# --- Remove the now-obsolete rows 37-40, shifting everything below up ---
$ws.Range("A37:F40").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp) | Out-Null
